# Update the hard-coded "datetime1" field text from 2020-09-09 to 2020-09-10
# across every slide layout and the slide master (the footer/date
# placeholders that carry a cached <a:fld type="datetime1"> value).

$p = $ppt.ActivePresentation

$oldDate = "2020-09-09"
$newDate = "2020-09-10"

function Update-DateFieldShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master's own date placeholder.
$master = $p.SlideMaster
Update-DateFieldShapes $master.Shapes

# Every custom (slide) layout hanging off the master.
$layouts = $master.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    $layout = $layouts.Item($j)
    Update-DateFieldShapes $layout.Shapes
}
